$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos: text replaced by the professor's name (content bug from the source edit) ---
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# --- Rows 13-25 get entirely rebuilt / renumbered, so wipe them first (Clear, not just ClearContents,
#     so cells that shouldn't exist afterwards don't leave empty stubs behind). ---
$ws.Range("A13:C25").Clear()

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "15/07/2016"
$ws.Range("C13").Value = "15/07/2016"

# Row 14
$ws.Range("A14").Value = "Short syllabus:"

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

# Row 16
$ws.Range("A16").Value = "Syllabus:"

# Row 17
$ws.Range("A17").Value = "Avaliação:"

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5817692 - Katia Cristiane Gandolpho Candioto"
$ws.Range("C18").Value = "5817692 - Katia Cristiane Gandolpho Candioto"

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "A avaliação será constituída por aulas expositivas, aulas de exercícios e práticas laboratoriais. Serão aplicadas pelo menos duas avaliações."
$ws.Range("C19").Value = "A avaliação será constituída por aulas expositivas, aulas de exercícios e práticas laboratoriais. Serão aplicadas pelo menos duas avaliações."

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final será a média das avaliações escritas e práticas"
$ws.Range("C20").Value = "A nota final será a média das avaliações escritas e práticas"

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A recuperação será uma prova escrita (RE) que comporá com a nota final (NF) a média final (MF), sendo MF = (NF + RE)/2."
$ws.Range("C21").Value = "A recuperação será uma prova escrita (RE) que comporá com a nota final (NF) a média final (MF), sendo MF = (NF + RE)/2."

# Row 22
$ws.Range("A22").Value = "Requisitos:"

# Row 23 (only B/C, no A)
$ws.Range("B23").Value = "LOM3011 -  Ensaios Mecânicos  (Requisito)`n"
$ws.Range("C23").Value = "LOM3011 -  Ensaios Mecânicos  (Requisito)`n"

# The sheet's <cols> definition has an overlapping range (min=1 max=2 style=1, then min=2 max=2
# style=2) and newly-created column-B cells inherit the wrong (first-matching) style. Column C
# and column A both pick up the correct style automatically, so just repair column B by pasting
# the number/cell format from a column-B cell that already carries the right style (style id 2).
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights (customHeight) for rows 13-23
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).EntireRow.AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).EntireRow.AutoFit()
$ws.Rows.Item(23).RowHeight = 30

# Rows 24-25 no longer exist (content now ends at row 23) - remove them so the sheet shrinks
# from 25 to 23 rows, shifting nothing else (they're already the last two rows).
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Delete()
